# "them tk hong nguyen thi" - add a new Facebook account (Hong Nguyen Thi)
# to the group-share list on Sheet1.
#
# Previously B1 held the plain-text label "Group_hongNguyenThi". It is
# replaced with the new account's e-mail address and turned into a
# mailto: hyperlink (matching the style already used by A1/C1), and the
# active selection moves to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new value for B1.
$ws.Range("B1").Value = "hongnguyen12229@gmail.com"

# Turn it into a mailto hyperlink, like A1 (copmapmap22@gmail.com) and
# C1 (hanhnhan7891@gmail.com) already are.
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:hongnguyen12229@gmail.com") | Out-Null

# Hyperlinks.Add resets formatting on the cell; re-apply the same
# "Hyperlink" look (centered, underlined, themed color) already used by
# the other header cells in row 1 by copying C1's format onto B1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null

# Move the active selection to B8, matching the edited workbook.
$ws.Range("B8").Select() | Out-Null
